# Applies the 05-11-2023 update to the england_premier-league_2023-2024 sheet.
# The scraped data for several fixtures in the same gameweek had been
# associated with the wrong rows (mismatched kickoff-time ordering); this
# script re-associates each row with its correct match data (teams, score,
# odds, odds timestamps and match URL) and appends the one fixture
# (Nottingham v Aston Villa) that was missing from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sets the "match" portion of a row (columns F..V) - home team, home goals,
# away team, away goals, and the 1x2 opening/closing odds with their
# timestamps, plus the match url. Columns A..E (index/pais/torneio/
# temporada/data_partida) are left untouched.
function SetMatchRow($row, $vals) {
    $cols = 6..22
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($row, $cols[$i]).Value2 = $vals[$i]
    }
}

# --- Re-pair rows with their correct match data ---------------------------

SetMatchRow 73 @("Everton", 3, "Bournemouth", 0, 1.69, "24/09/2023 10:02", 1.93, "07/10/2023 15:58", 4.01, "24/09/2023 10:02", 3.78, "07/10/2023 15:58", 4.57, "24/09/2023 10:02", 4.18, "07/10/2023 15:58", "https://www.betexplorer.com/football/england/premier-league/everton-bournemouth/CInUym42/")
SetMatchRow 75 @("Fulham", 3, "Sheffield Utd", 1, 1.61, "28/09/2023 14:26", 1.52, "07/10/2023 15:41", 3.93, "28/09/2023 14:26", 4.56, "07/10/2023 15:58", 5.39, "28/09/2023 14:26", 6.57, "07/10/2023 15:58", "https://www.betexplorer.com/football/england/premier-league/fulham-sheffield-utd/j9oYz7J8/")

SetMatchRow 83 @("Manchester City", 2, "Brighton", 1, 1.47, "01/10/2023 23:01", 1.36, "21/10/2023 15:34", 4.98, "01/10/2023 23:01", 5.69, "21/10/2023 15:58", 6.36, "01/10/2023 23:01", 7.74, "21/10/2023 15:58", "https://www.betexplorer.com/football/england/premier-league/manchester-city-brighton/ptI9zbPP/")
SetMatchRow 84 @("Bournemouth", 1, "Wolves", 2, 2.24, "01/10/2023 23:01", 2.31, "21/10/2023 15:55", 3.47, "01/10/2023 23:01", 3.46, "21/10/2023 15:57", 3.33, "01/10/2023 23:01", 3.32, "21/10/2023 15:58", "https://www.betexplorer.com/football/england/premier-league/bournemouth-wolves/bZIBFdm0/")
SetMatchRow 85 @("Brentford", 3, "Burnley", 0, 1.63, "02/10/2023 08:30", 1.74, "21/10/2023 15:45", 3.97, "02/10/2023 08:30", 3.9, "21/10/2023 15:58", 5.07, "02/10/2023 08:30", 5.09, "21/10/2023 15:58", "https://www.betexplorer.com/football/england/premier-league/brentford-burnley/6aMJDzIC/")
SetMatchRow 86 @("Nottingham", 2, "Luton", 2, 1.78, "02/10/2023 08:30", 1.76, "21/10/2023 15:56", 3.63, "02/10/2023 08:30", 3.65, "21/10/2023 15:58", 4.53, "02/10/2023 08:30", 5.39, "21/10/2023 15:58", "https://www.betexplorer.com/football/england/premier-league/nottingham-luton/tC3uVymm/")
SetMatchRow 87 @("Newcastle", 4, "Crystal Palace", 0, 1.49, "01/10/2023 23:01", 1.48, "21/10/2023 15:50", 4.51, "01/10/2023 23:01", 4.49, "21/10/2023 15:58", 7.02, "01/10/2023 23:01", 7.68, "21/10/2023 15:58", "https://www.betexplorer.com/football/england/premier-league/newcastle-utd-crystal-palace/2L4yWHXt/")

SetMatchRow 98 @("Liverpool", 3, "Nottingham", 0, 1.23, "10/10/2023 14:02", 1.22, "29/10/2023 14:51", 6.65, "10/10/2023 14:02", 7.31, "29/10/2023 14:59", 9.43, "10/10/2023 14:02", 12.08, "29/10/2023 14:59", "https://www.betexplorer.com/football/england/premier-league/liverpool-nottingham/IcEJreHn/")
SetMatchRow 99 @("Aston Villa", 3, "Luton", 1, 1.38, "10/10/2023 14:32", 1.34, "29/10/2023 14:55", 4.92, "10/10/2023 14:32", 5.85, "29/10/2023 14:55", 7.42, "10/10/2023 14:32", 8.55, "29/10/2023 14:57", "https://www.betexplorer.com/football/england/premier-league/aston-villa-luton/SSk1QD1I/")
SetMatchRow 100 @("Brighton", 1, "Fulham", 1, 1.51, "10/10/2023 14:02", 1.64, "29/10/2023 14:58", 4.73, "10/10/2023 14:02", 4.32, "29/10/2023 14:59", 6.2, "10/10/2023 14:02", 5.29, "29/10/2023 14:59", "https://www.betexplorer.com/football/england/premier-league/brighton-fulham/6Jl5PXGO/")

SetMatchRow 103 @("Sheffield Utd", 2, "Wolves", 1, 3.2, "23/10/2023 15:48", 4.35, "04/11/2023 15:59", 3.43, "23/10/2023 15:48", 3.76, "04/11/2023 15:59", 2.22, "23/10/2023 15:48", 1.87, "04/11/2023 15:59", "https://www.betexplorer.com/football/england/premier-league/sheffield-utd-wolves/0tW9gCV4/")
SetMatchRow 104 @("Brentford", 3, "West Ham", 2, 1.95, "21/10/2023 20:02", 2.16, "04/11/2023 15:50", 3.65, "21/10/2023 20:02", 3.71, "04/11/2023 15:50", 4, "21/10/2023 20:02", 3.39, "04/11/2023 15:50", "https://www.betexplorer.com/football/england/premier-league/brentford-west-ham/MkBzuDGB/")
SetMatchRow 106 @("Everton", 1, "Brighton", 1, 2.86, "21/10/2023 20:02", 2.86, "04/11/2023 15:59", 3.84, "21/10/2023 20:02", 3.52, "04/11/2023 15:59", 2.24, "21/10/2023 20:02", 2.54, "04/11/2023 15:59", "https://www.betexplorer.com/football/england/premier-league/everton-brighton/f39rwioO/")
SetMatchRow 107 @("Manchester City", 6, "Bournemouth", 1, 1.15, "21/10/2023 22:01", 1.09, "04/11/2023 15:50", 9.18, "21/10/2023 22:01", 11.5, "04/11/2023 15:21", 17.1, "21/10/2023 22:01", 28.5, "04/11/2023 15:52", "https://www.betexplorer.com/football/england/premier-league/manchester-city-bournemouth/AiwcdEon/")

# --- Append the new fixture (row 109) --------------------------------------

# Copy A108/E108 formatting (bold/border/center style for the index column,
# date/time number format for the match-date column) onto the new row before
# filling in the values.
$ws.Cells.Item(108, 1).Copy()
$ws.Cells.Item(109, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(108, 5).Copy()
$ws.Cells.Item(109, 5).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(109, 1).Value2 = 108
$ws.Cells.Item(109, 2).Value2 = "england"
$ws.Cells.Item(109, 3).Value2 = "premier-league"
$ws.Cells.Item(109, 4).Value2 = "2023-2024"
$ws.Cells.Item(109, 5).Value2 = 45235.625

SetMatchRow 109 @("Nottingham", 2, "Aston Villa", 0, 3.14, "21/10/2023 20:02", 3.93, "05/11/2023 14:35", 3.45, "21/10/2023 20:02", 3.81, "05/11/2023 14:35", 2.35, "21/10/2023 20:02", 1.95, "05/11/2023 14:35", "https://www.betexplorer.com/football/england/premier-league/nottingham-aston-villa/IHt5fhGb/")

$wb.Save()
